# Apply cryptocurrency price/volume updates scraped on Thu Jun  8 08:33:56 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.434.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.40%  "
$ws.Range("D3").Value = "'1.841.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'261.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.61%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.5213"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.32%  "
$ws.Range("E8").Value = "  -4.45%  "
$ws.Range("D9").Value = "'0.06786"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.36%  "
$ws.Range("D10").Value = "'18.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.02%  "
$ws.Range("D11").Value = "'0.7710"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.14%  "
$ws.Range("D12").Value = "'0.07718"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "'1.821.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.94%  "
$ws.Range("D14").Value = "'87.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.75%  "
$ws.Range("D15").Value = "'5.001"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.57%  "
$ws.Range("D16").Value = "'0.9996"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("E17").Value = "  -4.55%  "
$ws.Range("D19").Value = "'0.000007964"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("D20").Value = "'26.451.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("D21").Value = "'2.067.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("D22").Value = "'4.586"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.41%  "
$ws.Range("D23").Value = "'9.496"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.39%  "
$ws.Range("D24").Value = "'5.980"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.12%  "
$ws.Range("D25").Value = "'144.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("D26").Value = "'2.180"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.94%  "
$ws.Range("D27").Value = "'1.650"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").Value = "'16.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.11%  "
$ws.Range("D29").Value = "'111.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.01%  "
$ws.Range("D30").Value = "'4.203"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.50%  "
$ws.Range("D31").Value = "'4.125"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.18%  "
$ws.Range("D32").Value = "'0.08694"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.43%  "
$ws.Range("D33").Value = "'0.04796"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.45%  "
$ws.Range("D34").Value = "'1.129"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.76%  "
$ws.Range("D35").Value = "'0.7184"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.16%  "
$ws.Range("D36").Value = "'2.846"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("D37").Value = "'3.081"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.91%  "
$ws.Range("D38").Value = "'0.01781"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.58%  "
$ws.Range("D39").Value = "'2.213"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.46%  "
$ws.Range("D40").Value = "'0.4828"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.90%  "
$ws.Range("D41").Value = "'112.29"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.55%  "
$ws.Range("D42").Value = "'0.8996"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.94%  "
$ws.Range("D43").Value = "'6.068"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "'7.721"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.00%  "
$ws.Range("D46").Value = "'0.05890"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.4127"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.63%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "'35.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.21%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.949"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.17%  "
$ws.Range("D50").Value = "'0.1217"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.09%  "
$ws.Range("D51").Value = "'0.8865"
$ws.Range("D51").Style = "Normal"
